# The target revision re-distributes the full contents of rows 4, 6, 7, 8, 9,
# 10, 11 and 12 (columns A:AY) among themselves: every target row ends up
# holding exactly the data that used to live in a different, specific source
# row (a permutation of whole rows), rather than having individual cell
# values tweaked. This also carries along row-specific quirks such as the
# blank K/L/M/N marker cells and the "hack" public-comment cell that
# currently sit on row 7 - they must travel together with the rest of that
# row's data to its new home.
#
# Because this is a single closed permutation cycle, every source row has to
# be read out before any target row is overwritten. We therefore copy each
# affected row first into a scratch/staging area safely out of the way, and
# only then copy the staged rows into their final destinations, which avoids
# clobbering data that is still needed.
#
# Plain Range.Value2 assignment cannot be used for this, because assigning an
# empty string / $null through Value2 always deletes the cell instead of
# leaving a present-but-blank cell behind, and it cannot create a new blank
# cell either. Range.Copy, on the other hand, faithfully reproduces a
# present-but-empty cell as well as a genuinely absent one, so it is used for
# every single cell that is copied below. Presence of a cell is recorded
# up-front (by inspecting only the pristine, not-yet-touched source rows)
# since Value2 read back right after a Copy() is not reliable.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = 1   # column A
$lastCol  = 51  # column AY

# target row -> source row (i.e. the target row receives the values that
# currently sit in the source row)
$mapping = @{
    4  = 10
    6  = 11
    7  = 8
    8  = 9
    9  = 12
    10 = 7
    11 = 4
    12 = 6
}

# Pick a block of scratch rows well outside the populated area to stage the
# source rows' contents in.
$stagingBase = 2000
$stagingRowFor = @{}
$i = 0
foreach ($srcRow in $mapping.Values) {
    if (-not $stagingRowFor.ContainsKey($srcRow)) {
        $stagingRowFor[$srcRow] = $stagingBase + $i
        $i++
    }
}

# Step 0: record, for every distinct source row and every column, whether a
# cell is actually present there (Value2 -ne $null) *before* anything is
# copied anywhere. This presence map is reused later for the staging rows
# too, since a row staged from a source row has exactly the same shape.
$presence = @{}
foreach ($srcRow in $stagingRowFor.Keys) {
    $presence[$srcRow] = @{}
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $cell = $ws.Cells.Item($srcRow, $col)
        $presence[$srcRow][$col] = ($cell.Value2 -ne $null)
    }
}

# Step 1: copy every distinct source row, cell by cell, into its staging row.
# Only cells that were recorded as present are copied, so we never fabricate
# cells that were not present in the source row.
foreach ($srcRow in $stagingRowFor.Keys) {
    $stageRow = $stagingRowFor[$srcRow]
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        if ($presence[$srcRow][$col]) {
            $srcCell = $ws.Cells.Item($srcRow, $col)
            $dstCell = $ws.Cells.Item($stageRow, $col)
            $srcCell.Copy($dstCell)
        }
    }
}

# Step 2: clear the original target rows so that no stale cell from the
# "before" state lingers (e.g. row 7 must lose its K/L/M/N/AC cells once its
# old content has moved to row 10).
foreach ($targetRow in $mapping.Keys) {
    $rowRange = $ws.Range($ws.Cells.Item($targetRow, $firstCol), $ws.Cells.Item($targetRow, $lastCol))
    $rowRange.ClearContents()
}

# Step 3: copy from the staging rows into the real target rows, again relying
# on the presence map gathered in Step 0 (a staged row has the same shape as
# the original source row it came from).
foreach ($targetRow in $mapping.Keys) {
    $srcRow = $mapping[$targetRow]
    $stageRow = $stagingRowFor[$srcRow]
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        if ($presence[$srcRow][$col]) {
            $stageCell = $ws.Cells.Item($stageRow, $col)
            $dstCell = $ws.Cells.Item($targetRow, $col)
            $stageCell.Copy($dstCell)
        }
    }
}

# Step 4: remove the staging rows entirely so the worksheet's used range goes
# back to its original extent.
$stagingRows = $stagingRowFor.Values | Sort-Object -Descending
foreach ($stageRow in $stagingRows) {
    $ws.Rows.Item($stageRow).Delete() | Out-Null
}
